$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
# Paragraph "Write a function to validate and ensure ..." (Part 1 list)
# gets a darkMagenta highlight applied across the whole paragraph
# (all runs + the paragraph mark itself).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Write a function to validate and ensure*") {
        $p.Range.Font.HighlightColorIndex = 12   # wdDarkMagenta
    }
}

# --- Change 2 -------------------------------------------------------
# Paragraph "Aside from validating the user input in Part 1 ..."
# (Part 2 "General Things to Consider" list) also gets a darkMagenta
# highlight across the whole paragraph, and the "_GoBack" bookmark is
# moved from the following paragraph to the very start of this one.

# Remove the existing (empty) "_GoBack" bookmark wherever it currently
# sits.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Aside from validating the user input in Part 1*") {
        $p.Range.Font.HighlightColorIndex = 12   # wdDarkMagenta

        # Re-create the bookmark, collapsed, at the start of this paragraph.
        $startRange = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $startRange)
    }
}
